$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87
$prev = $row - 1

# Copy the previous row's cell so the new row inherits the same
# formatting / data typing (e.g. the date-formatted style on column A).
$ws.Range("A$prev`:H$prev").Copy($ws.Range("A$row`:H$row"))

$ws.Cells.Item($row, 1).Value = 45448.6262615741
$ws.Cells.Item($row, 2).Value = 9812
$ws.Cells.Item($row, 3).Value = 0.725000023841858
$ws.Cells.Item($row, 4).Value = 0.709999978542328
$ws.Cells.Item($row, 5).Value = 0.725000023841858
$ws.Cells.Item($row, 6).Value = 0.709999978542328

# Column G (adj_close) stores this number as text/shared-string in the
# source data, so copy it as a value from an existing cell holding the
# exact same text ("0.709999978542328") to avoid Excel auto-converting
# the literal into a numeric cell.
$ws.Cells.Item(85, 7).Copy()
$ws.Cells.Item($row, 7).PasteSpecial(-4163)

$ws.Cells.Item($row, 8).Value = "BWZ.MI"

$excel.CutCopyMode = 0
